# Updates leve profit calculation results across all job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed
# market-board pricing data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15236.412
$ws.Range("I21").Value = 16002.111
$ws.Range("J21").Value = 14375
$ws.Range("K21").Value = 16002.111
$ws.Range("L21").Value = 14375
$ws.Range("M21").Value = -15534.111
$ws.Range("N21").Value = -15311

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 15236.412
$ws.Range("I23").Value = 16002.111
$ws.Range("J23").Value = 14375
$ws.Range("K23").Value = 16002.111
$ws.Range("L23").Value = 14375
$ws.Range("M23").Value = -15768.111
$ws.Range("N23").Value = -14843

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 498.875
$ws.Range("I29").Value = 115.166664
$ws.Range("J29").Value = 1650
$ws.Range("K29").Value = 345.499992
$ws.Range("L29").Value = 4950
$ws.Range("M29").Value = -64.49999200000002
$ws.Range("N29").Value = -5512

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 443.33334
$ws.Range("I38").Value = 443.33334
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1330.00002
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -958.0000199999999
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 418.8125
$ws.Range("I39").Value = 131
$ws.Range("J39").Value = 1052
$ws.Range("K39").Value = 393
$ws.Range("L39").Value = 3156
$ws.Range("M39").Value = -97
$ws.Range("N39").Value = -3748

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1942.2222
$ws.Range("I40").Value = 1672.6666
$ws.Range("J40").Value = 2279.1667
$ws.Range("K40").Value = 1672.6666
$ws.Range("L40").Value = 2279.1667
$ws.Range("M40").Value = -1497.6666
$ws.Range("N40").Value = -2629.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4771.5
$ws.Range("I51").Value = 1885.8572
$ws.Range("J51").Value = 7657.143
$ws.Range("K51").Value = 1885.8572
$ws.Range("L51").Value = 7657.143
$ws.Range("M51").Value = -1401.8572
$ws.Range("N51").Value = -8625.143

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1347.9688
$ws.Range("I58").Value = 1013.5714
$ws.Range("J58").Value = 1986.3636
$ws.Range("K58").Value = 3040.7142
$ws.Range("L58").Value = 5959.0908
$ws.Range("M58").Value = -2890.7142
$ws.Range("N58").Value = -6259.0908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 66669304
$ws.Range("I131").Value = 100001900
$ws.Range("J131").Value = 4100
$ws.Range("K131").Value = 300005700
$ws.Range("L131").Value = 12300
$ws.Range("M131").Value = -300000660
$ws.Range("N131").Value = -22380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 559.6667
$ws.Range("I4").Value = 580
$ws.Range("J4").Value = 549.5
$ws.Range("K4").Value = 580
$ws.Range("L4").Value = 549.5
$ws.Range("M4").Value = -464
$ws.Range("N4").Value = -781.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 60008
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 60008
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 20833.334
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 28545.455
$ws.Range("J37").Value = 28545.455
$ws.Range("L37").Value = 28545.455
$ws.Range("N37").Value = -29091.455

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 31250
$ws.Range("J44").Value = 31250
$ws.Range("L44").Value = 31250
$ws.Range("N44").Value = -32226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 17415.143
$ws.Range("J55").Value = 15021.2
$ws.Range("L55").Value = 15021.2
$ws.Range("N55").Value = -15651.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 430.55554
$ws.Range("I22").Value = 397.5
$ws.Range("J22").Value = 695
$ws.Range("K22").Value = 397.5
$ws.Range("L22").Value = 695
$ws.Range("M22").Value = -224.5
$ws.Range("N22").Value = -1041

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 27499.75
$ws.Range("J74").Value = 27499.75
$ws.Range("L74").Value = 27499.75
$ws.Range("N74").Value = -29247.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 27499.75
$ws.Range("J77").Value = 27499.75
$ws.Range("L77").Value = 82499.25
$ws.Range("N77").Value = -91235.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 496.54544
$ws.Range("I34").Value = 356.6
$ws.Range("J34").Value = 613.1667
$ws.Range("K34").Value = 1069.8
$ws.Range("L34").Value = 1839.5001
$ws.Range("M34").Value = -985.8000000000002
$ws.Range("N34").Value = -2007.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3090
$ws.Range("J39").Value = 3090
$ws.Range("L39").Value = 9270
$ws.Range("N39").Value = -9858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 770.9091
$ws.Range("I55").Value = 450
$ws.Range("J55").Value = 786.1905
$ws.Range("K55").Value = 1350
$ws.Range("L55").Value = 2358.5715
$ws.Range("M55").Value = -1173
$ws.Range("N55").Value = -2712.5715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 29023
$ws.Range("J46").Value = 29023
$ws.Range("L46").Value = 29023
$ws.Range("N46").Value = -29335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 11539.8
$ws.Range("J57").Value = 11999.75
$ws.Range("L57").Value = 11999.75
$ws.Range("N57").Value = -13639.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 46205.047
$ws.Range("J22").Value = 786.2381
$ws.Range("L22").Value = 786.2381
$ws.Range("N22").Value = -1376.2381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 46205.047
$ws.Range("J27").Value = 786.2381
$ws.Range("L27").Value = 786.2381
$ws.Range("N27").Value = -1000.2381

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 835.26086
$ws.Range("I46").Value = 955.6667
$ws.Range("J46").Value = 757.8570999999999
$ws.Range("K46").Value = 955.6667
$ws.Range("L46").Value = 757.8570999999999
$ws.Range("M46").Value = -767.6667
$ws.Range("N46").Value = -1133.8571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 19283.533
$ws.Range("J104").Value = 19283.533
$ws.Range("L104").Value = 19283.533
$ws.Range("N104").Value = -26271.533

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 51548.75
$ws.Range("I137").Value = 58945
$ws.Range("J137").Value = 49083.332
$ws.Range("K137").Value = 58945
$ws.Range("L137").Value = 49083.332
$ws.Range("M137").Value = -53845
$ws.Range("N137").Value = -59283.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 77316.664
$ws.Range("I54").Value = 40000
$ws.Range("J54").Value = 84780
$ws.Range("K54").Value = 40000
$ws.Range("L54").Value = 84780
$ws.Range("M54").Value = -39480
$ws.Range("N54").Value = -85820
